$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row (16) that continues the chapter log, copying formatting
# from row 15 and filling in the new data.
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "3114-246"
$ws.Cells.Item(16, 3).Value = "Котел"
$ws.Cells.Item(16, 4).Value = "Мегу, Зак, Эван"
$ws.Cells.Item(16, 5).Value = "Эван зовет Мегу поужинать, вместо этого девушка, после разговора с Заком, решается открыть флешку Мэтта"

# Row 16 holds a long wrapped description, like the taller rows above it
# (e.g. row 12 at 150pt): match the row height Excel would auto-calculate
# for the wrapped text.
$ws.Rows.Item(16).RowHeight = 60

$ws.Range("E17").Select() | Out-Null
